# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with the latest scraped quotes. Price cells are stored as text
# (e.g. "28.469.03", "1.003") in the source sheet, so values that would
# otherwise be auto-parsed as numbers are written with a leading
# apostrophe to force Excel's text/quote-prefix interpretation and keep
# them as literal strings, matching the existing cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.469.03"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "'1.915.76"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'316.00"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.5111"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D8").Value = "'0.3966"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "'0.09707"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").Value = "'1.143"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").Value = "'42.14"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").Value = "'6.471"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'20.99"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'1.920.16"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "'7.398"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'0.00001131"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'93.92"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'0.06679"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'18.08"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'6.272"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").Value = "'28.528.26"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'11.45"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "'2.318"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").Value = "'2.666"
$ws.Range("E26").Value = "  +4.85%  "
$ws.Range("D27").Value = "'2.138.10"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D28").Value = "'21.19"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").Value = "'158.46"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "'128.43"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "'1.101"
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'5.690"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "'3.640"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "'9.803"
$ws.Range("E35").Value = "  +4.04%  "
$ws.Range("D36").Value = "'0.06704"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("D37").Value = "'0.02437"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").Value = "'1.257"
$ws.Range("E38").Value = "  +3.63%  "
$ws.Range("D39").Value = "'0.2219"
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("D40").Value = "'11.65"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").Value = "'0.6433"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").Value = "'5.055"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").Value = "'1.208"
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("D44").Value = "'1.003"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "'13.54"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "'0.6072"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "'3.782"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("D48").Value = "'1.283"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "'2.058"
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").Value = "'124.80"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'1.197"
$ws.Range("E51").Value = "  -0.32%  "
